$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) The "fixed" date shown in the Date placeholder on the slide master and
#    on every slide layout changes from 4/1/2016 to 4/4/2016.
# ---------------------------------------------------------------------------
function Update-FixedDate($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "4/1/2016") {
                $tr.Text = "4/4/2016"
            }
        }
    }
}

Update-FixedDate $p.SlideMaster.Shapes
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    Update-FixedDate $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) References slide - insert a new hyperlinked crunchify.com reference
#    (before the old one) and turn the former plain-text crunchify line into
#    a new plain-text misko.hevery.com reference.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(28)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$lastIndex = $tr.Paragraphs().Count
$crunchifyPara = $tr.Paragraphs($lastIndex, 1)

$crunchifyUrl = "http://crunchify.com/better-understanding-on-checked-vs-unchecked-exceptions-how-to-handle-exception-better-way-in-java"

# Split into two runs the same way the source deck does it elsewhere
# (main link text, then the trailing slash as its own run).
$newPara = $crunchifyPara.InsertBefore($crunchifyUrl + "`r")
$newPara.InsertAfter("/") | Out-Null

# Turn the whole new paragraph into a hyperlink pointing at the crunchify URL.
$hyperlinkPara = $tr.Paragraphs($lastIndex, 1)
$hyperlinkPara.ActionSettings(1).Hyperlink.Address = $crunchifyUrl + "/"

# The paragraph that used to hold the crunchify link now gets the
# misko.hevery.com reference as plain text (no hyperlink).
$miskoPara = $tr.Paragraphs($lastIndex + 1, 1)
$miskoChars = $miskoPara.Characters(1, $miskoPara.Length)
$miskoChars.Text = "http://misko.hevery.com/code-reviewers-guide/"
